$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 loses its yellow highlight fill entirely (A2:D2 were all shaded).
$ws.Range("A2:D2").ClearFormats()

# --- Values / formulas -------------------------------------------------
# B2 used to be "=B3"; it now holds a plain literal value (same number).
$ws.Range("B2").Value2 = 7699

# B3 / B4: the monthly inpatient cost input numbers were updated.
$ws.Range("B3").Value2 = 12577
$ws.Range("B4").Value2 = 12577

# C2:C4 now compute the daily->monthly conversion from the B column.
$ws.Range("C2").Formula = "=B2*((12358-6551)/3.92)/9454"
$ws.Range("C3").Formula = "=B3*((12358-6551)/3.92)/9454"
$ws.Range("C4").Formula = "=B4*((12358-6551)/3.92)/9454"

# D2 referenced a placeholder note; it now cites the same source as D3/D4.
$ws.Range("D2").Value2 = "skinner2018healthcare"

# --- Formatting ----------------------------------------------------------
# B2:B4 take on the "#,##0" thousands style (same look previously only on B3/B4).
$ws.Range("B2:B4").NumberFormat = "#,##0"
$ws.Range("B2:B4").WrapText = $true
$ws.Range("B2:B4").VerticalAlignment = -4160

# C2:C4 get an accounting-style comma format (0 decimals).
$ws.Range("C2:C4").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# D2:D4 share the wrap/top alignment already used by D3/D4.
$ws.Range("D2:D4").WrapText = $true
$ws.Range("D2:D4").VerticalAlignment = -4160

# Row heights for the data rows grew slightly.
$ws.Range("A2:D4").RowHeight = 17

# A blank, similarly-formatted helper cell appears just below the table.
$ws.Range("C5").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# --- View state ------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true) | Out-Null
$ws.Range("C4").Select() | Out-Null
